$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I20").Value = 0.4954102784692063
$ws.Range("J20").Value = 0.2573349478400102
$ws.Range("K20").Value = 0.276846350178154
$ws.Range("L20").Value = 2.488042276717461
